# Update column C ("Förändrad") date serial values from 45186 to 45188
# for all data rows on the active worksheet ("Avverkningsanmälningar").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
